$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header
$ws.Range("B1").Value = "Functionalization"

# Format column B (rows 2-5) as Text so numeric-looking values stay text
$ws.Range("B2:B5").NumberFormat = "@"

# Update row 2
$ws.Range("A2").Value = "skeleton_1"
$ws.Range("B2").Value = "1"
$ws.Range("C2").Value = "CH3"
$ws.Range("D2").Value = "CH3"
$ws.Range("E2").Value = "CH3"
$ws.Range("F2").Value = "CH3"

# Update row 3
$ws.Range("A3").Value = "skeleton_1"
$ws.Range("B3").Value = "2"
$ws.Range("C3").Value = "CCH3CH3CH3"
$ws.Range("D3").Value = "CCH3CH3CH3"
$ws.Range("E3").Value = "CCH3CH3CH3"
$ws.Range("F3").Value = "CCH3CH3CH3"

# New row 4
$ws.Range("A4").Value = "skeleton_2"
$ws.Range("B4").Value = "1"
$ws.Range("C4").Value = "CH3"
$ws.Range("D4").Value = "CH3"
$ws.Range("E4").Value = "CH3"
$ws.Range("F4").Value = "CH3"

# New row 5
$ws.Range("A5").Value = "skeleton_2"
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = "CCH3CH3CH3"
$ws.Range("D5").Value = "CCH3CH3CH3"
$ws.Range("E5").Value = "CCH3CH3CH3"
$ws.Range("F5").Value = "CCH3CH3CH3"

# Update dimension / ignored errors range to include new rows
$ws.Range("A1:F5").Cells.Item(1,1).Select() | Out-Null
